$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update participant name in first ceremony group (row 5)
$ws.Range("D5").Value = "Yaneth Mejía Rendon"

# 2. New ceremony row group (rows 8-10), mirroring the "Actividades Sprint 1"
#    group in rows 6-7 (copy formats down, then set values)
$ws.Range("A6:E7").Copy()
$ws.Range("A8:E10").PasteSpecial(-4122)

$ws.Range("A8").Value = 3
$ws.Range("B8").Value = 44314
$ws.Range("C8").Value = "Actividades Sprint 1"
$ws.Range("D8").Value = "Julio Edwin Mora Ardila"
$ws.Range("D9").Value = "Juan Carlos Rojas Buitrago"
$ws.Range("D10").Value = "Yaneth Mejía Rendon"

$ws.Range("A8:A10").Merge()
$ws.Range("B8:B10").Merge()
$ws.Range("C8:C10").Merge()
$ws.Range("E8:E10").Merge()

$ws.Hyperlinks.Add($ws.Range("E8"), "https://drive.google.com/file/d/1pKTFI8EGziGN2_UBOYa2Fwl9k36OQIIc/view?usp=sharing", "", "", "https://drive.google.com/file/d/1pKTFI8EGziGN2_UBOYa2Fwl9k36OQIIc/view?usp=sharing")
